$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLC Tags")

# New PLC tag rows to append under the existing table (columns: Name, Path,
# Data Type, Logical Address, Comment, Hmi Visible, Hmi Accessible,
# Hmi Writeable, Typeobject ID, Version ID).
$rows = @(
    @{ Name = "RPM_VAR1"; Addr = "%IW70" },
    @{ Name = "RPM_VAR2"; Addr = "%IW72" },
    @{ Name = "RPM_VAR3"; Addr = "%IW74" },
    @{ Name = "RPM_VAR4"; Addr = "%IW76" }
)

$startRow = 21
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.Name      # A - Name
    $ws.Cells.Item($r, 2).Value = "Entradas"     # B - Path
    $ws.Cells.Item($r, 3).Value = "Word"         # C - Data Type
    $ws.Cells.Item($r, 4).Value = $row.Addr      # D - Logical Address

    # E/I/J (Comment / Typeobject ID / Version ID) are blank text cells in
    # the source data, not "no cell at all" - force text typing with a
    # leading apostrophe, then drop the quote-prefix style it introduces so
    # the cell keeps the sheet's default style.
    foreach ($col in 5, 9, 10) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.Value = "'"
        $cell.Style = "Normal"
    }

    # F/G/H (Hmi Visible / Hmi Accessible / Hmi Writeable) hold the literal
    # text "True" - not the boolean TRUE. Assigning the string directly
    # gets auto-coerced to a Boolean cell by the COM layer, so instead
    # compute it as a text formula result and paste back as a value; this
    # keeps the cell's type as text without touching the stylesheet.
    foreach ($col in 6, 7, 8) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.Formula = '=""&"True"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = $false
